# BOM.xlsx update — "Add files via upload / added label maker"
#
# The TOOL LIST section (rows 78-80) already lists a few optional tools but
# had one blank, unused spacer row right before the AUTHORS block. This
# change turns that spacer row into a new tool-list entry ("Label maker"),
# inserts a fresh blank spacer row in its place (so the layout/spacing below
# is preserved), and bumps the revision footer from "Rev 1.0 / last update:
# 2023-02-20" to "Rev 1.01 / last update: 2023-02-21".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 82 (pushes the AUTHORS block and everything below
# it down by one row, duplicating the formatting of the row above — exactly
# what the previously-blank row 81 already had).
$ws.Rows(82).Insert()

# Row 81 (still blank after the insert) becomes the new tool-list item.
$ws.Range("A81").Value = "Label maker (optional but strongly recommended)"

# The revision/date footer (now at rows 87-88 after the shift) gets updated.
$ws.Range("A87").Value = "last update: 2023-02-21"
$ws.Range("A88").Value = "Rev 1.01"
